$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 527
$ws.Range("F5").Value = 2331
$ws.Range("F7").Value = 8193
$ws.Range("F8").Value = 121
$ws.Range("F10").Value = 1626
$ws.Range("F11").Value = 1335
$ws.Range("F12").Value = 214
$ws.Range("F13").Value = 4585
$ws.Range("F14").Value = 6192
$ws.Range("F15").Value = 813
$ws.Range("F16").Value = 63
$ws.Range("F17").Value = 1282
$ws.Range("F18").Value = 1282
$ws.Range("F19").Value = 491
$ws.Range("F20").Value = 6559
$ws.Range("F21").Value = 365
$ws.Range("F24").Value = 4419
$ws.Range("F25").Value = 330
$ws.Range("F26").Value = 729
$ws.Range("F27").Value = 2083
$ws.Range("F28").Value = 1207
$ws.Range("F29").Value = 363
$ws.Range("F30").Value = 1091
$ws.Range("F31").Value = 90
$ws.Range("F32").Value = 59
$ws.Range("F33").Value = 51
$ws.Range("F34").Value = 98
$ws.Range("F35").Value = 338
$ws.Range("F36").Value = 1205
$ws.Range("F37").Value = 1919
$ws.Range("F38").Value = 155
$ws.Range("F40").Value = 173
$ws.Range("F41").Value = 1244
$ws.Range("F42").Value = 567
$ws.Range("F43").Value = 75
$ws.Range("F44").Value = 1216
$ws.Range("F47").Value = 205
$ws.Range("F48").Value = 37
$ws.Range("F49").Value = 28

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 421
$ws.Range("F10").Value = 227
$ws.Range("F12").Value = 9
$ws.Range("F15").Value = 205
$ws.Range("F17").Value = 121
$ws.Range("F18").Value = 22
$ws.Range("F20").Value = 112
$ws.Range("F26").Value = 195

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 472
$ws.Range("F6").Value = 1598
$ws.Range("F7").Value = 499
$ws.Range("F8").Value = 3167
$ws.Range("F9").Value = 1108
$ws.Range("F10").Value = 1172
$ws.Range("F11").Value = 1555
$ws.Range("F12").Value = 1887
$ws.Range("F13").Value = 376
$ws.Range("F14").Value = 250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1598
$ws.Range("F5").Value = 527
$ws.Range("F6").Value = 499
$ws.Range("F7").Value = 3167
$ws.Range("F8").Value = 2331
$ws.Range("F9").Value = 121
$ws.Range("F10").Value = 1108
$ws.Range("F12").Value = 1626
$ws.Range("F13").Value = 1555
$ws.Range("F14").Value = 1335
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 1887
$ws.Range("F18").Value = 4585
$ws.Range("F20").Value = 421
$ws.Range("F21").Value = 813
$ws.Range("F22").Value = 63
$ws.Range("F23").Value = 1282
$ws.Range("F24").Value = 1282
$ws.Range("F25").Value = 491
$ws.Range("F26").Value = 6559
$ws.Range("F27").Value = 365
$ws.Range("F28").Value = 250
$ws.Range("F30").Value = 330
$ws.Range("F31").Value = 2083
$ws.Range("F32").Value = 1207
$ws.Range("F33").Value = 363
$ws.Range("F34").Value = 90
$ws.Range("F35").Value = 51
$ws.Range("F36").Value = 205
$ws.Range("F38").Value = 98
$ws.Range("F39").Value = 338
$ws.Range("F40").Value = 1919
$ws.Range("F41").Value = 155
$ws.Range("F43").Value = 1244
$ws.Range("F45").Value = 567
$ws.Range("F47").Value = 1216
$ws.Range("F49").Value = 205
